$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12
$ws.Range("F5").Value = 202
$ws.Range("F6").Value = 23
$ws.Range("F7").Value = 1047
$ws.Range("F8").Value = 829
$ws.Range("F9").Value = 241
$ws.Range("F12").Value = 833
$ws.Range("F13").Value = 285
$ws.Range("F15").Value = 503
$ws.Range("F18").Value = 1043
$ws.Range("F19").Value = 1193
$ws.Range("F20").Value = 2878
$ws.Range("F21").Value = 1430
$ws.Range("F22").Value = 704
$ws.Range("F23").Value = 197
$ws.Range("F24").Value = 1274
$ws.Range("F26").Value = 1020
$ws.Range("F27").Value = 360
$ws.Range("F28").Value = 3113
$ws.Range("F29").Value = 606
$ws.Range("F30").Value = 536
$ws.Range("F31").Value = 1405

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = "不可售"
$ws.Range("F4").Value = 367
$ws.Range("F5").Value = 14
$ws.Range("F11").Value = 25

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 741

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 12
$ws.Range("F3").Value = 741
$ws.Range("G4").Value = "不可售"
$ws.Range("F7").Value = 367
$ws.Range("F8").Value = 14
$ws.Range("F11").Value = 202
$ws.Range("F13").Value = 23
$ws.Range("F14").Value = 1047
$ws.Range("F15").Value = 829
$ws.Range("F16").Value = 241
$ws.Range("F21").Value = 25
$ws.Range("F24").Value = 833
$ws.Range("F25").Value = 285
$ws.Range("F27").Value = 503
$ws.Range("F30").Value = 1043
$ws.Range("F31").Value = 1193
$ws.Range("F32").Value = 2878
$ws.Range("F33").Value = 1430
$ws.Range("F34").Value = 704
$ws.Range("F35").Value = 197
$ws.Range("F36").Value = 1274
$ws.Range("F40").Value = 1020
$ws.Range("F41").Value = 360
$ws.Range("F42").Value = 3113
$ws.Range("F43").Value = 606
$ws.Range("F44").Value = 536
$ws.Range("F45").Value = 1405
